$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Replace all occurrences of "SARS-CoV-2" with "2019-nCoV" throughout the
#    document (there are four of them).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.Execute("SARS-CoV-2", $true, $false, $false, $false, $false, $true, 1, $false, "2019-nCoV", 2)

# ---------------------------------------------------------------------------
# 2) Remove the old "_Hlk34986703" bookmark that used to wrap the very first
#    "SARS-CoV-2" occurrence (title line). Removing it also causes the
#    remaining bookmarks further down in the document to be renumbered,
#    exactly like Word does.
# ---------------------------------------------------------------------------
$d.Bookmarks.Item("_Hlk34986703").Delete()

# ---------------------------------------------------------------------------
# 3) Merge the two runs "不排除感染，" and "可能会受到方法学局限性以及采样时机和"
#    (which used to be separated by the old "_GoBack" bookmark) into a single
#    run, while leaving the neighbouring runs ("：", "采集部位", "等的影响。")
#    untouched. We do this by temporarily bracketing the runs that must stay
#    untouched with throw-away bookmarks (acting as merge barriers), editing
#    across the old "_GoBack" location, and then removing the helper
#    bookmarks again (bookmark add/delete alone never triggers a run merge).
# ---------------------------------------------------------------------------
$findBarrier1 = $d.Content.Find
$findBarrier1.Execute("采集部位", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$barrierRange1 = $d.Range($findBarrier1.Parent.Start, $findBarrier1.Parent.Start)
$d.Bookmarks.Add("_TempBarrier1", $barrierRange1)

$findBarrier2 = $d.Content.Find
$findBarrier2.Execute("等的影响。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$barrierRange2 = $d.Range($findBarrier2.Parent.Start, $findBarrier2.Parent.Start)
$d.Bookmarks.Add("_TempBarrier2", $barrierRange2)

$d.Bookmarks.Item("_GoBack").Delete()

$findMerge = $d.Content.Find
$findMerge.Execute("可能会受到方法学局限性以及采样时机和", $true, $false, $false, $false, $false, $true, 1, $false, "可能会受到方法学局限性以及采样时机和", 2)

$d.Bookmarks.Item("_TempBarrier1").Delete()
$d.Bookmarks.Item("_TempBarrier2").Delete()

# ---------------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark at the location of the most recent
#    edit: right after the fourth "2019-nCoV" occurrence.
# ---------------------------------------------------------------------------
$searchRange = $d.Content
$searchRange.Start = 0
$occurrence = 0
while ($true) {
    $f = $searchRange.Find
    $found = $f.Execute("2019-nCoV", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }
    $occurrence += 1
    $lastEnd = $searchRange.End
    if ($occurrence -eq 4) { break }
    $searchRange.Start = $searchRange.End
    $searchRange.End = $d.Content.End
}

$goBackRange = $d.Range($lastEnd, $lastEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)
